# GBDS NOVEMBER FILES 2025 - fliqlo@GBDS
# Fill in the purchase entries for rows 19-23 on the "PE, NOVEMBER" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PE, NOVEMBER")

# The date column (C) is formatted as a date (m/d/yy) but the author always
# enters the dates as literal text (matching rows 11-18 already in the
# sheet). Flip to a text format while typing the value, then restore the
# original date format code so the cell keeps its original style/format.
foreach ($r in 19..23) {
    $ws.Range("C$r").NumberFormat = "@"
}

$ws.Range("C19").Value = "11/29/2025"
$ws.Range("C20").Value = "11/30/2025"
$ws.Range("C21").Value = "11/30/2025"
$ws.Range("C22").Value = "11/30/2025"
$ws.Range("C23").Value = "11/30/2025"

foreach ($r in 19..23) {
    $ws.Range("C$r").NumberFormat = "m/d/yy"
}

# Invoice numbers (column G)
$ws.Range("G19").Value = 518149982
$ws.Range("G20").Value = 518152453
$ws.Range("G21").Value = 518153334
$ws.Range("G22").Value = 518153364
$ws.Range("G23").Value = 518153400

# Gross purchase amounts (column I) - net of VAT credit memo deductions
$ws.Range("I19").Formula = "=1300225-56321.18"
$ws.Range("I20").Formula = "=122800-10224"
$ws.Range("I21").Formula = "=1353132-56255.04"
$ws.Range("I22").Formula = "=1353132-56255.04"
$ws.Range("I23").Formula = "=1595054-66248.54"

# Leave the selection where the author left it when saving.
$ws.Range("C24").Select()

$wb.Save()
